# Commit: "updated control signal names"
#
# The sheet's row 3 lists control-signal example values for the "addi"
# instruction. Column E (header "ALUControl (i_ALU_C)" in row 2) was left
# blank; this change fills it in with the placeholder value "x" (the ALU
# control value is a don't-care / not shown as a concrete signal here,
# mirrored by similar cells elsewhere in the row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "x"

# Leave the cursor where the author's last save shows it (row 3 pane,
# bottom-right frozen pane active cell).
$ws.Range("N6").Select()
